$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.104.30'
$ws.Range('D3').Value = '1.822.97'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4620'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3638'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.75%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07295'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8697'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.08'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.07%  '
$ws.Range('D12').Value = '1.876.69'
$ws.Range('E12').Value = '  +2.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07573'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.99%  '
$ws.Range('E14').Value = '  -2.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.33'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.484'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.60%  '
$ws.Range('E17').Value = '  -0.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008624'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.010'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').Value = '27.416.43'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('E21').Value = '  -2.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.206'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.36%  '
$ws.Range('E23').Value = '  -1.56%  '
$ws.Range('D24').Value = '2.094.21'
$ws.Range('E24').Value = '  +1.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.875'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.088'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '116.24'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.077'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08903'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.69%  '
$ws.Range('E32').Value = '  +0.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7326'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.457'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.136'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.61%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.473'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.87%  '
$ws.Range('E38').Value = '  -2.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05252'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.99%  '
$ws.Range('E40').Value = '  -2.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.931'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.140'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.98%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5195'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1630'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.37%  '
$ws.Range('E45').Value = '  -3.95%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4865'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.34%  '
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.15'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '103.51'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('E50').Value = '  -3.09%  '
$ws.Range('E51').Value = '  -1.05%  '
